$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force text storage via a temporary text format, then clear the format so
# the cell style reverts to its original (unstyled) state.
$textCells = @(
    'D4'
    'D5'
    'D6'
    'D8'
    'D11'
    'D14'
    'D19'
    'D20'
    'D21'
    'D22'
    'D24'
    'D28'
    'D29'
    'D31'
    'D35'
    'D39'
    'D40'
    'D41'
    'D44'
    'D45'
    'D46'
    'D48'
    'D50'
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply new values
$ws.Range('D2').Value = '69.282.62'
$ws.Range('E2').Value = '  -2.82%  '
$ws.Range('D3').Value = '3.683.27'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '684.68'
$ws.Range('E5').Value = '  -3.05%  '
$ws.Range('D6').Value = '162.29'
$ws.Range('E6').Value = '  -5.57%  '
$ws.Range('D7').Value = '3.682.04'
$ws.Range('E7').Value = '  -3.52%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E10').Value = '  -8.51%  '
$ws.Range('D11').Value = '7.36'
$ws.Range('E11').Value = '  -3.68%  '
$ws.Range('E12').Value = '  -3.84%  '
$ws.Range('E13').Value = '  -5.41%  '
$ws.Range('D14').Value = '33.67'
$ws.Range('E14').Value = '  -6.62%  '
$ws.Range('D15').Value = '4.302.65'
$ws.Range('E15').Value = '  -3.56%  '
$ws.Range('D16').Value = '3.682.84'
$ws.Range('E16').Value = '  -3.24%  '
$ws.Range('D17').Value = '69.350.91'
$ws.Range('E17').Value = '  -2.66%  '
$ws.Range('E18').Value = '  -1.16%  '
$ws.Range('D19').Value = '16.38'
$ws.Range('E19').Value = '  -6.47%  '
$ws.Range('D20').Value = '6.63'
$ws.Range('D21').Value = '482.79'
$ws.Range('E21').Value = '  -6.86%  '
$ws.Range('D22').Value = '9.95'
$ws.Range('E22').Value = '  -6.92%  '
$ws.Range('D24').Value = '80.27'
$ws.Range('E24').Value = '  -5.26%  '
$ws.Range('D25').Value = '3.827.87'
$ws.Range('E25').Value = '  -3.53%  '
$ws.Range('E26').Value = '  -9.88%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').Value = '11.43'
$ws.Range('E28').Value = '  -5.39%  '
$ws.Range('D29').Value = '9.52'
$ws.Range('E29').Value = '  -8.75%  '
$ws.Range('E30').Value = '  -10.32%  '
$ws.Range('D31').Value = '2.72'
$ws.Range('E31').Value = '  -10.37%  '
$ws.Range('E33').Value = '  -7.57%  '
$ws.Range('E34').Value = '  -2.59%  '
$ws.Range('D35').Value = '27.13'
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('D37').Value = '3.652.41'
$ws.Range('E37').Value = '  -3.57%  '
$ws.Range('E38').Value = '  -7.65%  '
$ws.Range('D39').Value = '6.29'
$ws.Range('E39').Value = '  +5.34%  '
$ws.Range('D40').Value = '2.32'
$ws.Range('E40').Value = '  -2.15%  '
$ws.Range('D41').Value = '0.0936'
$ws.Range('E41').Value = '  -7.92%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('D44').Value = '0.950'
$ws.Range('E44').Value = '  -7.19%  '
$ws.Range('D45').Value = '163.13'
$ws.Range('E45').Value = '  -2.49%  '
$ws.Range('D46').Value = '48.38'
$ws.Range('E46').Value = '  -1.76%  '
$ws.Range('E47').Value = '  -13.63%  '
$ws.Range('D48').Value = '30.00'
$ws.Range('E48').Value = '  +4.63%  '
$ws.Range('E49').Value = '  -8.09%  '
$ws.Range('D50').Value = '1.33'
$ws.Range('E50').Value = '  -1.61%  '
$ws.Range('E51').Value = '  -3.45%  '

# Restore original (default) formatting on the forced-text cells
foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
